$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force a literal text value into the cell without altering its
    # NumberFormat/style: build a formula that evaluates to the exact
    # string, then paste-special as values only so the formula is
    # discarded and only the literal text remains.
    $escaped = $text.Replace('"', '""')
    $range.Formula = '=""&"' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$ws.Range('D2').Value = '43.177.90'
$ws.Range('E2').Value = '  +2.63%  '
$ws.Range('D3').Value = '2.295.57'
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '252.72'
$ws.Range('E5').Value = '  +0.77%  '
Set-TextValue $ws.Range('D6') '0.640'
$ws.Range('E6').Value = '  +2.52%  '
Set-TextValue $ws.Range('D7') '74.12'
$ws.Range('E7').Value = '  +8.74%  '
$ws.Range('E8').Value = '  +0.00%  '
Set-TextValue $ws.Range('D9') '0.647'
$ws.Range('E9').Value = '  +4.28%  '
Set-TextValue $ws.Range('D10') '39.28'
$ws.Range('E10').Value = '  -0.57%  '
Set-TextValue $ws.Range('D11') '0.0985'
$ws.Range('E11').Value = '  +4.75%  '
Set-TextValue $ws.Range('D12') '59.07'
$ws.Range('E12').Value = '  -0.44%  '
Set-TextValue $ws.Range('D13') '7.44'
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('D15').Value = '2.639.75'
$ws.Range('E15').Value = '  +3.43%  '
Set-TextValue $ws.Range('D16') '15.32'
$ws.Range('E16').Value = '  +5.61%  '
Set-TextValue $ws.Range('D17') '0.876'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '2.294.28'
$ws.Range('E18').Value = '  +3.65%  '
$ws.Range('D19').Value = '43.064.77'
$ws.Range('E19').Value = '  +2.52%  '
Set-TextValue $ws.Range('D20') '0.0000101'
$ws.Range('E20').Value = '  +4.65%  '
Set-TextValue $ws.Range('D21') '6.33'
$ws.Range('E21').Value = '  +3.11%  '
Set-TextValue $ws.Range('D22') '72.63'
$ws.Range('E22').Value = '  +0.29%  '
Set-TextValue $ws.Range('D23') '237.97'
$ws.Range('E23').Value = '  +2.63%  '
$ws.Range('E24').Value = '  +8.96%  '
$ws.Range('E25').Value = '  +0.46%  '
Set-TextValue $ws.Range('D26') '11.59'
$ws.Range('E27').Value = '  -0.10%  '
Set-TextValue $ws.Range('D28') '2.43'
$ws.Range('E28').Value = '  +1.13%  '
Set-TextValue $ws.Range('D29') '3.65'
$ws.Range('E29').Value = '  -0.75%  '
Set-TextValue $ws.Range('D30') '2.20'
$ws.Range('E30').Value = '  -0.15%  '
Set-TextValue $ws.Range('D31') '167.12'
$ws.Range('E31').Value = '  +0.17%  '
Set-TextValue $ws.Range('D32') '21.10'
$ws.Range('E32').Value = '  +2.96%  '
$ws.Range('E33').Value = '  +5.78%  '
Set-TextValue $ws.Range('D34') '0.129'
$ws.Range('E34').Value = '  +5.84%  '
Set-TextValue $ws.Range('D35') '0.0818'
$ws.Range('E35').Value = '  +4.20%  '
Set-TextValue $ws.Range('D36') '31.45'
$ws.Range('E36').Value = '  +18.27%  '
Set-TextValue $ws.Range('D37') '0.127'
$ws.Range('E37').Value = '  +3.91%  '
$ws.Range('E38').Value = '  +12.33%  '
Set-TextValue $ws.Range('D39') '4.77'
$ws.Range('E39').Value = '  +3.57%  '
Set-TextValue $ws.Range('D40') '0.0310'
$ws.Range('E40').Value = '  -1.80%  '
Set-TextValue $ws.Range('D41') '14.52'
$ws.Range('E41').Value = '  +19.26%  '
$ws.Range('E42').Value = '  +4.88%  '
Set-TextValue $ws.Range('D43') '5.97'
$ws.Range('E43').Value = '  +4.78%  '
Set-TextValue $ws.Range('D44') '0.219'
$ws.Range('E44').Value = '  +11.76%  '
Set-TextValue $ws.Range('D45') '62.10'
$ws.Range('E45').Value = '  +0.00%  '
Set-TextValue $ws.Range('D46') '9.11'
$ws.Range('E46').Value = '  +6.07%  '
Set-TextValue $ws.Range('D47') '4.93'
$ws.Range('E47').Value = '  -3.06%  '
$ws.Range('E48').Value = '  +3.49%  '
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('E50').Value = '  +2.16%  '
Set-TextValue $ws.Range('D51') '98.77'
$ws.Range('E51').Value = '  +6.00%  '

$excel.CutCopyMode = 0
